$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Diff: cell C10 value changes from 18 to 1 (stored as numeric 1.0, t="n")
$ws.Range("C10").Value = 1
